$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 372.1085406666666
$ws.Cells.Item(2, 8).Value = 1116.325622
$ws.Cells.Item(2, 9).Value = 0.8095247142929753
$ws.Cells.Item(2, 10).Value = 0.8095247142929753
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 229.5846506666667
$ws.Cells.Item(2, 14).Value = 688.753952
$ws.Cells.Item(2, 15).Value = 0.5033187360873315
$ws.Cells.Item(2, 16).Value = 0.5033187360873315
$ws.Cells.Item(2, 17).Value = 85430.40931903978
$ws.Cells.Item(2, 18).Value = 768873.6838713581
$ws.Cells.Item(2, 19).Value = 0.4074489560293985
$ws.Cells.Item(2, 20).Value = 0.4074489560293985

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 372.1085406666666
$ws.Cells.Item(3, 8).Value = 1116.325622
$ws.Cells.Item(3, 9).Value = 0.8095247142929753
$ws.Cells.Item(3, 10).Value = 0.8095247142929753
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 135.7283196666666
$ws.Cells.Item(3, 14).Value = 407.1849589999999
$ws.Cells.Item(3, 15).Value = 0.2975573763642838
$ws.Cells.Item(3, 16).Value = 0.2975573763642838
$ws.Cells.Item(3, 17).Value = 50505.66695830215
$ws.Cells.Item(3, 18).Value = 454551.0026247193
$ws.Cells.Item(3, 19).Value = 0.2408800500870641
$ws.Cells.Item(3, 20).Value = 0.2408800500870641

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 372.1085406666666
$ws.Cells.Item(4, 8).Value = 1116.325622
$ws.Cells.Item(4, 9).Value = 0.8095247142929753
$ws.Cells.Item(4, 10).Value = 0.8095247142929753
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 90.23148833333335
$ws.Cells.Item(4, 14).Value = 270.694465
$ws.Cells.Item(4, 15).Value = 0.1978146123067711
$ws.Cells.Item(4, 16).Value = 0.1978146123067711
$ws.Cells.Item(4, 17).Value = 33575.90744589802
$ws.Cells.Item(4, 18).Value = 302183.1670130822
$ws.Cells.Item(4, 19).Value = 0.1601358175106145
$ws.Cells.Item(4, 20).Value = 0.1601358175106145

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 372.1085406666666
$ws.Cells.Item(5, 8).Value = 1116.325622
$ws.Cells.Item(5, 9).Value = 0.8095247142929753
$ws.Cells.Item(5, 10).Value = 0.8095247142929753
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.5972149999999999
$ws.Cells.Item(5, 14).Value = 1.791645
$ws.Cells.Item(5, 15).Value = 0.001309275241613694
$ws.Cells.Item(5, 16).Value = 0.001309275241613694
$ws.Cells.Item(5, 17).Value = 222.2288021142433
$ws.Cells.Item(5, 18).Value = 2000.05921902819
$ws.Cells.Item(5, 19).Value = 0.001059890665898192
$ws.Cells.Item(5, 20).Value = 0.001059890665898192

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 69.70494733333334
$ws.Cells.Item(6, 8).Value = 209.114842
$ws.Cells.Item(6, 9).Value = 0.1516435969830949
$ws.Cells.Item(6, 10).Value = 0.1516435969830949
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 229.5846506666667
$ws.Cells.Item(6, 14).Value = 688.753952
$ws.Cells.Item(6, 15).Value = 0.5033187360873315
$ws.Cells.Item(6, 16).Value = 0.5033187360873315
$ws.Cells.Item(6, 17).Value = 16003.18598326173
$ws.Cells.Item(6, 18).Value = 144028.6738493556
$ws.Cells.Item(6, 19).Value = 0.07632506356926798
$ws.Cells.Item(6, 20).Value = 0.07632506356926798

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 69.70494733333334
$ws.Cells.Item(7, 8).Value = 209.114842
$ws.Cells.Item(7, 9).Value = 0.1516435969830949
$ws.Cells.Item(7, 10).Value = 0.1516435969830949
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 135.7283196666666
$ws.Cells.Item(7, 14).Value = 407.1849589999999
$ws.Cells.Item(7, 15).Value = 0.2975573763642838
$ws.Cells.Item(7, 16).Value = 0.2975573763642838
$ws.Cells.Item(7, 17).Value = 9460.935374006829
$ws.Cells.Item(7, 18).Value = 85148.41836606147
$ws.Cells.Item(7, 19).Value = 0.04512267086073252
$ws.Cells.Item(7, 20).Value = 0.04512267086073252

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 69.70494733333334
$ws.Cells.Item(8, 8).Value = 209.114842
$ws.Cells.Item(8, 9).Value = 0.1516435969830949
$ws.Cells.Item(8, 10).Value = 0.1516435969830949
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 90.23148833333335
$ws.Cells.Item(8, 14).Value = 270.694465
$ws.Cells.Item(8, 15).Value = 0.1978146123067711
$ws.Cells.Item(8, 16).Value = 0.1978146123067711
$ws.Cells.Item(8, 17).Value = 6289.581142083282
$ws.Cells.Item(8, 18).Value = 56606.23027874954
$ws.Cells.Item(8, 19).Value = 0.02999731934601515
$ws.Cells.Item(8, 20).Value = 0.02999731934601515

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 69.70494733333334
$ws.Cells.Item(9, 8).Value = 209.114842
$ws.Cells.Item(9, 9).Value = 0.1516435969830949
$ws.Cells.Item(9, 10).Value = 0.1516435969830949
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.5972149999999999
$ws.Cells.Item(9, 14).Value = 1.791645
$ws.Cells.Item(9, 15).Value = 0.001309275241613694
$ws.Cells.Item(9, 16).Value = 0.001309275241613694
$ws.Cells.Item(9, 17).Value = 41.62884012167667
$ws.Cells.Item(9, 18).Value = 374.65956109509
$ws.Cells.Item(9, 19).Value = 0.0001985432070792113
$ws.Cells.Item(9, 20).Value = 0.0001985432070792113

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.387356
$ws.Cells.Item(10, 8).Value = 1.162068
$ws.Cells.Item(10, 9).Value = 0.0008426956679571845
$ws.Cells.Item(10, 10).Value = 0.0008426956679571845
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 229.5846506666667
$ws.Cells.Item(10, 14).Value = 688.753952
$ws.Cells.Item(10, 15).Value = 0.5033187360873315
$ws.Cells.Item(10, 16).Value = 0.5033187360873315
$ws.Cells.Item(10, 17).Value = 88.93099194363734
$ws.Cells.Item(10, 18).Value = 800.3789274927361
$ws.Cells.Item(10, 19).Value = 0.0004241445185024797
$ws.Cells.Item(10, 20).Value = 0.0004241445185024797

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.387356
$ws.Cells.Item(11, 8).Value = 1.162068
$ws.Cells.Item(11, 9).Value = 0.0008426956679571845
$ws.Cells.Item(11, 10).Value = 0.0008426956679571845
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 135.7283196666666
$ws.Cells.Item(11, 14).Value = 407.1849589999999
$ws.Cells.Item(11, 15).Value = 0.2975573763642838
$ws.Cells.Item(11, 16).Value = 0.2975573763642838
$ws.Cells.Item(11, 17).Value = 52.57517899280133
$ws.Cells.Item(11, 18).Value = 473.176610935212
$ws.Cells.Item(11, 19).Value = 0.0002507503120308875
$ws.Cells.Item(11, 20).Value = 0.0002507503120308875

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 0.387356
$ws.Cells.Item(12, 8).Value = 1.162068
$ws.Cells.Item(12, 9).Value = 0.0008426956679571845
$ws.Cells.Item(12, 10).Value = 0.0008426956679571845
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 90.23148833333335
$ws.Cells.Item(12, 14).Value = 270.694465
$ws.Cells.Item(12, 15).Value = 0.1978146123067711
$ws.Cells.Item(12, 16).Value = 0.1978146123067711
$ws.Cells.Item(12, 17).Value = 34.95170839484668
$ws.Cells.Item(12, 18).Value = 314.5653755536201
$ws.Cells.Item(12, 19).Value = 0.0001666975168495459
$ws.Cells.Item(12, 20).Value = 0.0001666975168495459

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 0.387356
$ws.Cells.Item(13, 8).Value = 1.162068
$ws.Cells.Item(13, 9).Value = 0.0008426956679571845
$ws.Cells.Item(13, 10).Value = 0.0008426956679571845
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 0.5972149999999999
$ws.Cells.Item(13, 14).Value = 1.791645
$ws.Cells.Item(13, 15).Value = 0.001309275241613694
$ws.Cells.Item(13, 16).Value = 0.001309275241613694
$ws.Cells.Item(13, 17).Value = 0.23133481354
$ws.Cells.Item(13, 18).Value = 2.08201332186
$ws.Cells.Item(13, 19).Value = 0.000001103320574271456
$ws.Cells.Item(13, 20).Value = 0.000001103320574271456

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 0.425022
$ws.Cells.Item(14, 8).Value = 1.275066
$ws.Cells.Item(14, 9).Value = 0.0009246383125251667
$ws.Cells.Item(14, 10).Value = 0.0009246383125251667
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 229.5846506666667
$ws.Cells.Item(14, 14).Value = 688.753952
$ws.Cells.Item(14, 15).Value = 0.5033187360873315
$ws.Cells.Item(14, 16).Value = 0.5033187360873315
$ws.Cells.Item(14, 17).Value = 97.57852739564801
$ws.Cells.Item(14, 18).Value = 878.206746560832
$ws.Cells.Item(14, 19).Value = 0.0004653877867980899
$ws.Cells.Item(14, 20).Value = 0.0004653877867980899

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 0.425022
$ws.Cells.Item(15, 8).Value = 1.275066
$ws.Cells.Item(15, 9).Value = 0.0009246383125251667
$ws.Cells.Item(15, 10).Value = 0.0009246383125251667
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 135.7283196666666
$ws.Cells.Item(15, 14).Value = 407.1849589999999
$ws.Cells.Item(15, 15).Value = 0.2975573763642838
$ws.Cells.Item(15, 16).Value = 0.2975573763642838
$ws.Cells.Item(15, 17).Value = 57.68752188136599
$ws.Cells.Item(15, 18).Value = 519.1876969322939
$ws.Cells.Item(15, 19).Value = 0.0002751329503608872
$ws.Cells.Item(15, 20).Value = 0.0002751329503608872

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 0.425022
$ws.Cells.Item(16, 8).Value = 1.275066
$ws.Cells.Item(16, 9).Value = 0.0009246383125251667
$ws.Cells.Item(16, 10).Value = 0.0009246383125251667
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 90.23148833333335
$ws.Cells.Item(16, 14).Value = 270.694465
$ws.Cells.Item(16, 15).Value = 0.1978146123067711
$ws.Cells.Item(16, 16).Value = 0.1978146123067711
$ws.Cells.Item(16, 17).Value = 38.35036763441001
$ws.Cells.Item(16, 18).Value = 345.1533087096901
$ws.Cells.Item(16, 19).Value = 0.0001829069693161529
$ws.Cells.Item(16, 20).Value = 0.0001829069693161529

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 0.425022
$ws.Cells.Item(17, 8).Value = 1.275066
$ws.Cells.Item(17, 9).Value = 0.0009246383125251667
$ws.Cells.Item(17, 10).Value = 0.0009246383125251667
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 0.5972149999999999
$ws.Cells.Item(17, 14).Value = 1.791645
$ws.Cells.Item(17, 15).Value = 0.001309275241613694
$ws.Cells.Item(17, 16).Value = 0.001309275241613694
$ws.Cells.Item(17, 17).Value = 0.25382951373
$ws.Cells.Item(17, 18).Value = 2.28446562357
$ws.Cells.Item(17, 19).Value = 0.000001210606050036666
$ws.Cells.Item(17, 20).Value = 0.000001210606050036666

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 17.03711166666666
$ws.Cells.Item(18, 8).Value = 51.111335
$ws.Cells.Item(18, 9).Value = 0.03706435474344739
$ws.Cells.Item(18, 10).Value = 0.03706435474344739
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 229.5846506666667
$ws.Cells.Item(18, 14).Value = 688.753952
$ws.Cells.Item(18, 15).Value = 0.5033187360873315
$ws.Cells.Item(18, 16).Value = 0.5033187360873315
$ws.Cells.Item(18, 17).Value = 3911.459330360658
$ws.Cells.Item(18, 18).Value = 35203.13397324592
$ws.Cells.Item(18, 19).Value = 0.01865518418336443
$ws.Cells.Item(18, 20).Value = 0.01865518418336443

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 17.03711166666666
$ws.Cells.Item(19, 8).Value = 51.111335
$ws.Cells.Item(19, 9).Value = 0.03706435474344739
$ws.Cells.Item(19, 10).Value = 0.03706435474344739
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 135.7283196666666
$ws.Cells.Item(19, 14).Value = 407.1849589999999
$ws.Cells.Item(19, 15).Value = 0.2975573763642838
$ws.Cells.Item(19, 16).Value = 0.2975573763642838
$ws.Cells.Item(19, 17).Value = 2312.418538490029
$ws.Cells.Item(19, 18).Value = 20811.76684641026
$ws.Cells.Item(19, 19).Value = 0.0110287721540953
$ws.Cells.Item(19, 20).Value = 0.0110287721540953

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 17.03711166666666
$ws.Cells.Item(20, 8).Value = 51.111335
$ws.Cells.Item(20, 9).Value = 0.03706435474344739
$ws.Cells.Item(20, 10).Value = 0.03706435474344739
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 13).Value = 90.23148833333335
$ws.Cells.Item(20, 14).Value = 270.694465
$ws.Cells.Item(20, 15).Value = 0.1978146123067711
$ws.Cells.Item(20, 16).Value = 0.1978146123067711
$ws.Cells.Item(20, 17).Value = 1537.28394258453
$ws.Cells.Item(20, 18).Value = 13835.55548326078
$ws.Cells.Item(20, 19).Value = 0.007331870963975676
$ws.Cells.Item(20, 20).Value = 0.007331870963975676

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 17.03711166666666
$ws.Cells.Item(21, 8).Value = 51.111335
$ws.Cells.Item(21, 9).Value = 0.03706435474344739
$ws.Cells.Item(21, 10).Value = 0.03706435474344739
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 13).Value = 0.5972149999999999
$ws.Cells.Item(21, 14).Value = 1.791645
$ws.Cells.Item(21, 15).Value = 0.001309275241613694
$ws.Cells.Item(21, 16).Value = 0.001309275241613694
$ws.Cells.Item(21, 17).Value = 10.17481864400833
$ws.Cells.Item(21, 18).Value = 91.57336779607499
$ws.Cells.Item(21, 19).Value = 0.00004852744201198276
$ws.Cells.Item(21, 20).Value = 0.00004852744201198276
